$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column D (Corequisites, Concurrent, Recommended),
# shifting the existing "Terms Typically Offered" column from D to G.
$ws.Range("D1:F1").EntireColumn.Insert()

# Set header text for new columns
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill "NA" for the new columns for each data row (2-14)
$ws.Range("D2:F14").Value = "NA"
